$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '52.605.89'
$ws.Range('E2').Value = '  -13.05%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.305.24'
$ws.Range('E3').Value = '  -20.54%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '443.36'
$ws.Range('E5').Value = '  -15.76%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '121.77'
$ws.Range('E6').Value = '  -14.20%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.997'
$ws.Range('E7').Value = '  -0.20%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.465'
$ws.Range('E8').Value = '  -15.35%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.311.01'
$ws.Range('E9').Value = '  -20.56%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.29'
$ws.Range('E10').Value = '  -11.53%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0874'
$ws.Range('E11').Value = '  -18.45%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.301'
$ws.Range('E12').Value = '  -16.24%  '
$ws.Range('E13').Value = '  -6.13%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '52.617.04'
$ws.Range('E14').Value = '  -13.03%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '18.77'
$ws.Range('E15').Value = '  -17.22%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000118'
$ws.Range('E16').Value = '  -15.86%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.328.01'
$ws.Range('E17').Value = '  -19.72%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.93'
$ws.Range('E18').Value = '  -21.02%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '296.65'
$ws.Range('E19').Value = '  -15.98%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '8.89'
$ws.Range('E20').Value = '  -23.60%  '
$ws.Range('E21').Value = '  -0.06%  '
$ws.Range('E22').Value = '  -1.74%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.12'
$ws.Range('E23').Value = '  -21.72%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '53.61'
$ws.Range('E24').Value = '  -16.90%  '
$ws.Range('E25').Value = '  -17.15%  '
$ws.Range('E26').Value = '  -19.83%  '
$ws.Range('B27').Value = 'InternetComputer(DFINITY)'
$ws.Range('C27').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.85'
$ws.Range('E27').Value = '  -12.56%  '
$ws.Range('B28').Value = 'USDe'
$ws.Range('C28').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.996'
$ws.Range('E28').Value = '  -0.35%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0₃0670'
$ws.Range('E29').Value = '  -20.23%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '141.35'
$ws.Range('E30').Value = '  -5.62%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '16.84'
$ws.Range('E31').Value = '  -14.04%  '
$ws.Range('E32').Value = '  -20.03%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.73'
$ws.Range('E33').Value = '  -15.26%  '
$ws.Range('E34').Value = '  -17.06%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.40'
$ws.Range('E35').Value = '  -21.32%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.996'
$ws.Range('E36').Value = '  -0.10%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.991'
$ws.Range('E37').Value = '  -17.40%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '31.64'
$ws.Range('E38').Value = '  -16.14%  '
$ws.Range('B39').Value = 'Mantle'
$ws.Range('C39').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.568'
$ws.Range('E39').Value = '  -12.27%  '
$ws.Range('B40').Value = 'WhiteBITCoin'
$ws.Range('C40').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '10.15'
$ws.Range('E40').Value = '  -1.65%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0502'
$ws.Range('E41').Value = '  -13.56%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.12'
$ws.Range('E42').Value = '  -15.91%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.898.17'
$ws.Range('E43').Value = '  -17.05%  '
$ws.Range('E44').Value = '  -19.58%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0205'
$ws.Range('E45').Value = '  -13.40%  '
$ws.Range('B46').Value = 'RenderToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.21'
$ws.Range('E46').Value = '  -14.73%  '
$ws.Range('B47').Value = 'Stellar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0820'
$ws.Range('E47').Value = '  -10.81%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '15.52'
$ws.Range('E48').Value = '  -23.79%  '
$ws.Range('E49').Value = '  -5.19%  '
$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '14.99'
$ws.Range('E50').Value = '  -17.94%  '
$ws.Range('B51').Value = 'Cosmos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.42'
$ws.Range('E51').Value = '  -13.67%  '
